$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.351.23'
$ws.Range("D3").Value = '3.570.76'
$ws.Range("E3").Value = '  +2.00%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '619.78'
$ws.Range("E5").Value = '  +2.71%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '155.58'
$ws.Range("E6").Value = '  +4.01%  '
$ws.Range("D7").Value = '3.568.06'
$ws.Range("E7").Value = '  +1.94%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  +2.25%  '
$ws.Range("E10").Value = '  +5.68%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.43'
$ws.Range("E11").Value = '  +7.12%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.439'
$ws.Range("E12").Value = '  +4.17%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '33.28'
$ws.Range("E13").Value = '  +5.77%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000222'
$ws.Range("E14").Value = '  +1.22%  '
$ws.Range("D15").Value = '4.169.78'
$ws.Range("E15").Value = '  +1.89%  '
$ws.Range("D16").Value = '3.563.92'
$ws.Range("E16").Value = '  +1.68%  '
$ws.Range("D17").Value = '68.371.48'
$ws.Range("E17").Value = '  +1.63%  '
$ws.Range("E18").Value = '  -0.08%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.74'
$ws.Range("E19").Value = '  +6.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.03'
$ws.Range("E20").Value = '  +6.96%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.09'
$ws.Range("E21").Value = '  +12.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '455.44'
$ws.Range("E22").Value = '  +2.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.645'
$ws.Range("E23").Value = '  +4.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.50'
$ws.Range("E24").Value = '  +1.52%  '
$ws.Range("E25").Value = '  +2.73%  '
$ws.Range("D26").Value = '3.708.94'
$ws.Range("E26").Value = '  +1.89%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  -0.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.28'
$ws.Range("E28").Value = '  +12.74%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.57'
$ws.Range("E29").Value = '  +4.40%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.72'
$ws.Range("E30").Value = '  +11.25%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.57'
$ws.Range("E31").Value = '  +3.62%  '
$ws.Range("E32").Value = '  +4.03%  '
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.40'
$ws.Range("E34").Value = '  +4.94%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '26.23'
$ws.Range("E35").Value = '  +2.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.93'
$ws.Range("E36").Value = '  +4.88%  '
$ws.Range("D37").Value = '3.559.97'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.26'
$ws.Range("E38").Value = '  +3.72%  '
$ws.Range("E39").Value = '  +9.12%  '
$ws.Range("E40").Value = '  -0.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '181.75'
$ws.Range("E41").Value = '  +5.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0920'
$ws.Range("E42").Value = '  +5.17%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("E43").Value = '  +0.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.61'
$ws.Range("E44").Value = '  +4.10%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '31.22'
$ws.Range("E45").Value = '  +12.76%  '
$ws.Range("E46").Value = '  +2.21%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '46.37'
$ws.Range("E47").Value = '  +2.11%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.34'
$ws.Range("E48").Value = '  +5.09%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.68'
$ws.Range("E49").Value = '  +4.94%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.82'
$ws.Range("E50").Value = '  +3.81%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.263'
$ws.Range("E51").Value = '  +8.15%  '
